$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(104, 6).Value = 77   # F104
$ws.Cells.Item(104, 7).Value = 5460.84   # G104
$ws.Cells.Item(118, 6).Value = 94   # F118
$ws.Cells.Item(118, 7).Value = 6610.08   # G118
$ws.Cells.Item(140, 6).Value = 32   # F140
$ws.Cells.Item(140, 7).Value = 1430.4   # G140
$ws.Cells.Item(143, 2).Value = 333209.02   # B143
$ws.Cells.Item(163, 6).Value = 14   # F163
$ws.Cells.Item(163, 7).Value = 558.04   # G163
$ws.Cells.Item(173, 6).Value = 6   # F173
$ws.Cells.Item(173, 7).Value = 146.58   # G173
$ws.Cells.Item(175, 6).Value = 15   # F175
$ws.Cells.Item(175, 7).Value = 289.2   # G175
$ws.Cells.Item(176, 2).Value = 23073.54   # B176
$ws.Cells.Item(180, 6).Value = 0   # F180
$ws.Cells.Item(180, 7).Value = 0   # G180
$ws.Cells.Item(186, 2).Value = 14128.79   # B186
$ws.Cells.Item(193, 6).Value = 106   # F193
$ws.Cells.Item(193, 7).Value = 5244.88   # G193
$ws.Cells.Item(201, 6).Value = 164   # F201
$ws.Cells.Item(201, 7).Value = 7301.28   # G201
$ws.Cells.Item(205, 2).Value = 34110.34   # B205
$ws.Cells.Item(284, 6).Value = 12   # F284
$ws.Cells.Item(284, 7).Value = 224.28   # G284
$ws.Cells.Item(288, 2).Value = 2941.8   # B288
$ws.Cells.Item(298, 6).Value = 126   # F298
$ws.Cells.Item(298, 7).Value = 10688.58   # G298
$ws.Cells.Item(303, 2).Value = 21465.19   # B303
$ws.Cells.Item(315, 6).Value = 57   # F315
$ws.Cells.Item(315, 7).Value = 8179.5   # G315
$ws.Cells.Item(320, 6).Value = 7   # F320
$ws.Cells.Item(320, 7).Value = 460.18   # G320
$ws.Cells.Item(323, 2).Value = 39984.11   # B323
$ws.Cells.Item(345, 6).Value = 5   # F345
$ws.Cells.Item(345, 7).Value = 302.9   # G345
$ws.Cells.Item(369, 6).Value = 27   # F369
$ws.Cells.Item(369, 7).Value = 6679.26   # G369
$ws.Cells.Item(374, 6).Value = 15   # F374
$ws.Cells.Item(374, 7).Value = 2033.1   # G374
$ws.Cells.Item(380, 6).Value = 57   # F380
$ws.Cells.Item(380, 7).Value = 1964.22   # G380
$ws.Cells.Item(386, 2).Value = 146774.92   # B386
$ws.Cells.Item(391, 2).Value = 57077   # B391
$ws.Cells.Item(391, 4).Value = 93.08   # D391
$ws.Cells.Item(391, 5).Value = 111.2   # E391
$ws.Cells.Item(391, 6).Value = 1   # F391
$ws.Cells.Item(391, 7).Value = 93.08   # G391
$ws.Cells.Item(392, 2).Value = 61610   # B392
$ws.Cells.Item(392, 4).Value = 102.71   # D392
$ws.Cells.Item(392, 5).Value = 122.71   # E392
$ws.Cells.Item(392, 6).Value = 360   # F392
$ws.Cells.Item(392, 7).Value = 36975.6   # G392
$ws.Cells.Item(400, 6).Value = 11   # F400
$ws.Cells.Item(400, 7).Value = 1745.37   # G400
$ws.Cells.Item(409, 6).Value = 168   # F409
$ws.Cells.Item(409, 7).Value = 23027.76   # G409
$ws.Cells.Item(414, 6).Value = 238   # F414
$ws.Cells.Item(414, 7).Value = 11338.32   # G414
$ws.Cells.Item(415, 6).Value = 169   # F415
$ws.Cells.Item(415, 7).Value = 8673.08   # G415
$ws.Cells.Item(424, 6).Value = 367   # F424
$ws.Cells.Item(424, 7).Value = 19131.71   # G424
$ws.Cells.Item(429, 6).Value = 60   # F429
$ws.Cells.Item(429, 7).Value = 5143.2   # G429
$ws.Cells.Item(434, 6).Value = 49   # F434
$ws.Cells.Item(434, 7).Value = 5100.41   # G434
$ws.Cells.Item(438, 6).Value = 248   # F438
$ws.Cells.Item(438, 7).Value = 25072.8   # G438
$ws.Cells.Item(442, 6).Value = 64   # F442
$ws.Cells.Item(442, 7).Value = 3486.72   # G442
$ws.Cells.Item(445, 6).Value = 110   # F445
$ws.Cells.Item(445, 7).Value = 13049.3   # G445
$ws.Cells.Item(446, 6).Value = 242   # F446
$ws.Cells.Item(446, 7).Value = 14309.46   # G446
$ws.Cells.Item(451, 6).Value = 536   # F451
$ws.Cells.Item(451, 7).Value = 11625.84   # G451
$ws.Cells.Item(452, 6).Value = 130   # F452
$ws.Cells.Item(452, 7).Value = 781.3   # G452
$ws.Cells.Item(460, 6).Value = 143   # F460
$ws.Cells.Item(460, 7).Value = 21617.31   # G460
$ws.Cells.Item(471, 2).Value = 589506.64   # B471
$ws.Cells.Item(531, 6).Value = 73   # F531
$ws.Cells.Item(531, 7).Value = 5011.45   # G531
$ws.Cells.Item(538, 6).Value = 154   # F538
$ws.Cells.Item(538, 7).Value = 4579.96   # G538
$ws.Cells.Item(546, 2).Value = 71651.61   # B546
$ws.Cells.Item(566, 6).Value = 0   # F566
$ws.Cells.Item(566, 7).Value = 0   # G566
$ws.Cells.Item(567, 6).Value = 91   # F567
$ws.Cells.Item(567, 7).Value = 16966.04   # G567
$ws.Cells.Item(575, 2).Value = 75179.35000000001   # B575
$ws.Cells.Item(584, 6).Value = 27   # F584
$ws.Cells.Item(584, 7).Value = 938.25   # G584
$ws.Cells.Item(586, 2).Value = 6826.3   # B586
$ws.Cells.Item(592, 6).Value = 592   # F592
$ws.Cells.Item(592, 7).Value = 7784.8   # G592
$ws.Cells.Item(593, 6).Value = 689   # F593
$ws.Cells.Item(593, 7).Value = 8826.09   # G593
$ws.Cells.Item(609, 2).Value = 121406.68   # B609
$ws.Cells.Item(678, 6).Value = 367   # F678
$ws.Cells.Item(678, 7).Value = 6066.51   # G678
$ws.Cells.Item(681, 2).Value = 48708.88   # B681
$ws.Cells.Item(731, 6).Value = 53   # F731
$ws.Cells.Item(731, 7).Value = 2116.82   # G731
$ws.Cells.Item(743, 2).Value = 13597.99   # B743
$ws.Cells.Item(746, 6).Value = 34   # F746
$ws.Cells.Item(746, 7).Value = 9719.58   # G746
$ws.Cells.Item(757, 6).Value = 144   # F757
$ws.Cells.Item(757, 7).Value = 3945.6   # G757
$ws.Cells.Item(760, 6).Value = 31   # F760
$ws.Cells.Item(760, 7).Value = 1815.36   # G760
$ws.Cells.Item(765, 2).Value = 86739.92999999999   # B765
$ws.Cells.Item(771, 6).Value = 162   # F771
$ws.Cells.Item(771, 7).Value = 4406.4   # G771
$ws.Cells.Item(773, 6).Value = 128   # F773
$ws.Cells.Item(773, 7).Value = 3481.6   # G773
$ws.Cells.Item(774, 2).Value = 81067.91   # B774
$ws.Cells.Item(800, 6).Value = 420   # F800
$ws.Cells.Item(800, 7).Value = 6661.2   # G800
$ws.Cells.Item(808, 2).Value = 51991.84   # B808
$ws.Cells.Item(842, 6).Value = 105   # F842
$ws.Cells.Item(842, 7).Value = 8439.9   # G842
$ws.Cells.Item(843, 2).Value = 27004.42   # B843
$ws.Cells.Item(850, 6).Value = 14   # F850
$ws.Cells.Item(850, 7).Value = 8186.08   # G850
$ws.Cells.Item(855, 2).Value = 29937.54   # B855
$ws.Cells.Item(868, 6).Value = 16   # F868
$ws.Cells.Item(868, 7).Value = 1068.64   # G868
$ws.Cells.Item(870, 2).Value = 8223.219999999999   # B870
$ws.Cells.Item(901, 6).Value = 91   # F901
$ws.Cells.Item(901, 7).Value = 10085.53   # G901
$ws.Cells.Item(905, 2).Value = 109406.11   # B905
$ws.Cells.Item(924, 6).Value = 301   # F924
$ws.Cells.Item(924, 7).Value = 9066.120000000001   # G924
$ws.Cells.Item(930, 6).Value = 111   # F930
$ws.Cells.Item(930, 7).Value = 4088.13   # G930
$ws.Cells.Item(937, 2).Value = 67358.46000000001   # B937
$ws.Cells.Item(943, 6).Value = 227   # F943
$ws.Cells.Item(943, 7).Value = 8489.799999999999   # G943
$ws.Cells.Item(946, 2).Value = 29419.85   # B946
$ws.Cells.Item(956, 6).Value = 4   # F956
$ws.Cells.Item(956, 7).Value = 22424.52   # G956
$ws.Cells.Item(961, 6).Value = 5   # F961
$ws.Cells.Item(961, 7).Value = 54267.55   # G961
$ws.Cells.Item(962, 2).Value = 191522.6   # B962
$ws.Cells.Item(1016, 2).Value = 4192411.77   # B1016
$ws.Cells.Item(1017, 2).Value = 4192411.77   # B1017
